$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.42580000000003
$ws.Range("A14").Value = -20.49659999999998
$ws.Range("A21").Value = -21.36170000000001
$ws.Range("A23").Value = -21.37150000000002
$ws.Range("A25").Value = -22.36730000000003
